# Update results.xlsx with the new SNOPT/IPOPT ("ampl") rows.
#
# In the original sheet, rows 23 and 24 are an empty gap, and the SNOPT /
# IPOPT summary rows live at 25 / 26. The new layout removes that gap
# (so SNOPT/IPOPT become rows 23/24), fills in the tv-norm (E) and
# obj-with-tv-norm (F) columns for both of them, relabels them with their
# own "Energy2_..." run names in column B, updates the runtime in G24, and
# extends the A-column merged label down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank 23:24 gap -> old rows 25/26 (SNOPT/IPOPT) shift up to 23/24.
$ws.Rows("23:24").Delete()

# Pick up the column A (merged-label) and column F (0.000 number format)
# cell formatting used by the rest of the table so the new rows match.
$ws.Range("A22").Copy()
$ws.Range("A23:A24").PasteSpecial(-4122)
$ws.Range("F22").Copy()
$ws.Range("F23:F24").PasteSpecial(-4122)

# Row 23: SNOPT continuous-relaxation result.
$ws.Range("B23").Value = "Energy2_evotime2.0_n_ts40_ptypeCONSTANT_offset0.5_snopt"
$ws.Range("E23").Value = 0.928
$ws.Range("F23").Formula = "=D23+E23*0.01"

# Row 24: IPOPT continuous-relaxation result.
$ws.Range("B24").Value = "Energy2_evotime2.0_n_ts40_ptypeCONSTANT_offset0.5_Ipopt"
$ws.Range("E24").Value = 0.928
$ws.Range("F24").Formula = "=D24+E24*0.01"
$ws.Range("G24").Value = 0.07

# Extend the merged label in column A from A2:A22 to A2:A24.
$ws.Range("A2:A22").MergeCells = $false
$ws.Range("A2:A24").Merge()

# Match the saved selection/view state.
$ws.Range("B11").Select() | Out-Null
